$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Update the "panel_query_time" timestamps on the "data" sheet
#    (column F, rows 2-18): 10:51:15.xxx -> 14:34:10.xxx
# ---------------------------------------------------------------------------
$dataSheet = $wb.Worksheets.Item("data")

$newTimes = @(
    "2021-10-05 14:34:10.095204",
    "2021-10-05 14:34:10.095213",
    "2021-10-05 14:34:10.095216",
    "2021-10-05 14:34:10.095219",
    "2021-10-05 14:34:10.095222",
    "2021-10-05 14:34:10.095224",
    "2021-10-05 14:34:10.095227",
    "2021-10-05 14:34:10.095230",
    "2021-10-05 14:34:10.095233",
    "2021-10-05 14:34:10.095235",
    "2021-10-05 14:34:10.095238",
    "2021-10-05 14:34:10.095241",
    "2021-10-05 14:34:10.095243",
    "2021-10-05 14:34:10.095246",
    "2021-10-05 14:34:10.095249",
    "2021-10-05 14:34:10.095252",
    "2021-10-05 14:34:10.095254"
)

for ($i = 0; $i -lt $newTimes.Length; $i++) {
    $row = $i + 2
    $dataSheet.Cells.Item($row, 6).Value = $newTimes[$i]
}

# ---------------------------------------------------------------------------
# 2. Add a new "metadata" worksheet right after "data"
# ---------------------------------------------------------------------------
$metaSheet = $wb.Worksheets.Add($null, $dataSheet)
$metaSheet.Name = "metadata"

# Header row
$metaSheet.Cells.Item(1, 2).Value = "data_name"
$metaSheet.Cells.Item(1, 3).Value = "data_id"
$metaSheet.Cells.Item(1, 4).Value = "data_version"
$metaSheet.Cells.Item(1, 5).Value = "data_version_created"
$metaSheet.Cells.Item(1, 6).Value = "panel_query_time"
$metaSheet.Cells.Item(1, 7).Value = "panel_get_request"

# Data row
$metaSheet.Cells.Item(2, 1).Value = 0
$metaSheet.Cells.Item(2, 2).Value = "Imprinting disorders"
$metaSheet.Cells.Item(2, 3).Value = 3663

# data_version must be stored as text "0.8" (not the number 0.8)
$metaSheet.Cells.Item(2, 4).NumberFormat = "@"
$metaSheet.Cells.Item(2, 4).Value = "0.8"
$metaSheet.Cells.Item(2, 4).ClearFormats()

$metaSheet.Cells.Item(2, 5).Value = "2021-09-22T06:50:28.236759Z"
$metaSheet.Cells.Item(2, 6).Value = "2021-10-05 14:34:10.091543"
$metaSheet.Cells.Item(2, 7).Value = "https://panelapp.agha.umccr.org/api/v1/panels/3663/?format=json"

# ---------------------------------------------------------------------------
# 3. Match the bold/centered/bordered header style used on "data" (B1:F1, A2)
#    by copying that formatting onto the new sheet's header row and A2 cell.
# ---------------------------------------------------------------------------
$dataSheet.Range("B1:F1").Copy() | Out-Null
$metaSheet.Range("B1:G1").PasteSpecial(-4122) | Out-Null
$dataSheet.Range("A2").Copy() | Out-Null
$metaSheet.Range("A2").PasteSpecial(-4122) | Out-Null

# ---------------------------------------------------------------------------
# 4. Keep "data" as the active sheet / selection, as it was originally.
# ---------------------------------------------------------------------------
$dataSheet.Activate()
$dataSheet.Range("A1").Select() | Out-Null
